# Add three new worksheets (strategy_id-6003, strategy_id-6004, strategy_id-6005)
# after the existing "strategy_id-0" sheet, each populated like the "General /
# frac_gnrl_eating_red_meat" row (row 6) of strategy_id-0, but with new values
# in the V:AS (year 12 onward) columns.

$wb = $excel.ActiveWorkbook
$base = $wb.Worksheets.Item("strategy_id-0")

$sheetNames = @("strategy_id-6003", "strategy_id-6004", "strategy_id-6005")

# V2:AS2 replacement values for each new sheet (24 values each, columns V..AS)
$valuesBySheet = @{
    "strategy_id-6003" = @(0.8475,0.8058333333333333,0.7649999999999999,0.7250000000000001,0.6858333333333333,0.6475,0.6099999999999999,0.5733333333333335,0.5375,0.5024999999999999,0.4683333333333334,0.435,0.4025000000000001,0.3708333333333333,0.34,0.31,0.2808333333333333,0.2525,0.225,0.2,0.175,0.15,0.125,0.1)
    "strategy_id-6004" = @(0.8558333333333333,0.8225,0.7899999999999999,0.7583333333333334,0.7274999999999999,0.6974999999999999,0.6683333333333332,0.6400000000000001,0.6125,0.5858333333333333,0.5600000000000001,0.535,0.5108333333333334,0.4875,0.465,0.4433333333333334,0.4225,0.4025,0.3833333333333333,0.3666666666666666,0.35,0.3333333333333333,0.3166666666666666,0.3)
    "strategy_id-6005" = @(0.8558333333333333,0.8225,0.7899999999999999,0.7583333333333334,0.7274999999999999,0.6974999999999999,0.6683333333333332,0.6400000000000001,0.6125,0.5858333333333333,0.5600000000000001,0.535,0.5108333333333334,0.4875,0.465,0.4433333333333334,0.4225,0.4025,0.3833333333333333,0.3666666666666666,0.35,0.3333333333333333,0.3166666666666666,0.3)
}

$cols = @("V","W","X","Y","Z","AA","AB","AC","AD","AE","AF","AG","AH","AI","AJ","AK","AL","AM","AN","AO","AP","AQ","AR","AS")

$prev = $base
foreach ($name in $sheetNames) {
    $newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $prev)
    $newSheet.Name = $name

    # Header row (row 1): copy values + formats from base sheet's header row.
    $base.Range("A1:AS1").Copy()
    $newSheet.Range("A1").PasteSpecial()
    $base.Range("A1:AS1").Copy()
    $newSheet.Range("A1").PasteSpecial(-4122)

    # Data row (row 2): copy values + formats from base sheet's row 6
    # ("General" / "frac_gnrl_eating_red_meat").
    $base.Range("A6:AS6").Copy()
    $newSheet.Range("A2").PasteSpecial()
    $base.Range("A6:AS6").Copy()
    $newSheet.Range("A2").PasteSpecial(-4122)

    # Base row 6 has no values in D6:G6 - remove the blank cells the paste
    # created so the sparse row matches the source layout.
    $newSheet.Range("D2:G2").Clear()

    # Overwrite the V2:AS2 values with the new per-sheet trajectory values.
    $vals = $valuesBySheet[$name]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $newSheet.Range($cols[$i] + "2").Value = $vals[$i]
    }

    $prev = $newSheet
}

# Keep the original sheet as the active/selected tab, matching the source.
$base.Activate()

Write-Host "done"
